# Append new scrape results (2025-10-23 18:26:46 JST) into the
# "ランサーズ" worksheet.
#
# Two brand-new listings were discovered by the scraper and were merged
# into the existing (score-sorted) list, landing at rows 4 and 6.
# Every other previously-known listing shifts down by the number of new
# rows inserted above it, and every row's "取得日時" (fetched-at) timestamp
# is refreshed to the new run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Insert the two new rows at their sorted positions. Doing this one at a
# time (in ascending order) pushes every following row down correctly.
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(6).Insert()

$timestamp = "2025-10-23 18:26:46"

# Full data set for rows 2-14 after the insert: A=取得日時, B=タイトル,
# C=カテゴリ, D=価格, E=締切, F=URL, G=優先度スコア, H=スキル概要
$rows = @(
    @(2,  "【26年5月/17日間/対面】Python Webアプリ開発 研修講師募集(カリキュラム設計含む)", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5419191", 295, "🔥Python ◆開発 ◇アプリ"),
    @(3,  "【急募】APIを利用した診断サイト構築のフリーランス募集", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5418643", 220, "🔥API ◇サイト"),
    @(4,  "【新規開発】iOS向けポイ活アプリの技術検証・設計者募集", "システム開発", "5,000,000 円 ~ / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5419221", 100, "◆開発 ◇アプリ"),
    @(5,  "【急募】施行主向け建築資材配達アプリ開発者を募集します", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5418447", 100, "◆開発 ◇アプリ"),
    @(6,  "ボディスコアアプリの開発", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5419226", 93, "◆開発 ◇アプリ"),
    @(7,  "【高品質な恋愛マッチングアプリ制作】エンジニア募集", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5418455", 45, "◇アプリ"),
    @(8,  "【Webarena suiteX/DNS】ドメイン設定変更によるウェブサイト分割とサイト切り替え", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5417544", 30, "◇サイト"),
    @(9,  "ERPシステムの第三者技術検証・品質評価報告書作成", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5418891", 40, $null),
    @(10, "【急募】セッション体験を再現するクローンシステム構築依頼", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5418644", 40, $null),
    @(11, "【急募】既存の予約システムの料金修正を依頼します", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5418759", 25, $null),
    @(12, "Stable Diffusion LoRA制作依頼 画風再現+キャラLoRA量産テンプレ構築", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5418738", 18, $null),
    @(13, "【メールマーケティング】戦略立案・実行者募集", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5418443", 18, $null),
    @(14, "【急募】HPの微修正をお手伝いしてくれる方募集!", "システム開発", "5,000 円 ~", "期限情報なし", "https://www.lancers.jp/work/detail/5418445", 10, $null)
)

foreach ($r in $rows) {
    $rowIndex = $r[0]
    $ws.Cells.Item($rowIndex, 1).Value = $timestamp
    $ws.Cells.Item($rowIndex, 2).Value = $r[1]
    $ws.Cells.Item($rowIndex, 3).Value = $r[2]
    $ws.Cells.Item($rowIndex, 4).Value = $r[3]
    $ws.Cells.Item($rowIndex, 5).Value = $r[4]
    $ws.Cells.Item($rowIndex, 6).Value = $r[5]
    $ws.Cells.Item($rowIndex, 7).Value = $r[6]
    if ($r[7] -ne $null) {
        $ws.Cells.Item($rowIndex, 8).Value = $r[7]
    }
}

# Rebuild the hyperlinks collection (row-insert does not shift hyperlink
# refs automatically in this runtime), in row order F2..F14.
$ws.Hyperlinks.Delete()
foreach ($r in $rows) {
    $rowIndex = $r[0]
    $ws.Hyperlinks.Add($ws.Cells.Item($rowIndex, 6), $r[5]) | Out-Null
}

# Hyperlinks.Add() registers a fresh (duplicate) cell style for the link
# font; restore the original shared "Hyperlink" style on the URL column
# so the cells keep using the workbook's existing style slot.
$ws.Range("F2:F14").Style = "Hyperlink"
